{"js": "// Update the date line and the 25 division-problem cells in the table.\n//\n// Every \"before\" value below is unique in the document, so each search\n// matches exactly one spot. To avoid any ambiguity from new values that\n// happen to equal other (still unprocessed) old values \u2014 e.g. 76\u00f73=\n// becomes 28\u00f75=, which is itself the *old* text of the very next cell \u2014\n// we resolve ALL search ranges against the original document first (one\n// batch of `search()` calls + a single `context.sync()`), and only then\n// issue the `insertText` replace calls. That way every Range object is\n// anchored to its original location before any text is mutated.\n\nconst replacements = [\n  [\"2025-09-18 Thursday\", \"2025-09-19 Friday\"],\n  [\"87\\u00F77=\", \"32\\u00F74=\"],\n  [\"88\\u00F75=\", \"61\\u00F75=\"],\n  [\"60\\u00F76=\", \"23\\u00F77=\"],\n  [\"79\\u00F73=\", \"80\\u00F75=\"],\n  [\"67\\u00F74=\", \"84\\u00F74=\"],\n  [\"93\\u00F78=\", \"74\\u00F73=\"],\n  [\"34\\u00F77=\", \"38\\u00F79=\"],\n  [\"38\\u00F77=\", \"87\\u00F78=\"],\n  [\"15\\u00F72=\", \"46\\u00F76=\"],\n  [\"35\\u00F76=\", \"20\\u00F75=\"],\n  [\"91\\u00F76=\", \"17\\u00F79=\"],\n  [\"30\\u00F78=\", \"16\\u00F74=\"],\n  [\"99\\u00F72=\", \"69\\u00F74=\"],\n  [\"80\\u00F77=\", \"37\\u00F72=\"],\n  [\"14\\u00F78=\", \"21\\u00F74=\"],\n  [\"76\\u00F73=\", \"28\\u00F75=\"],\n  [\"28\\u00F75=\", \"64\\u00F75=\"],\n  [\"23\\u00F75=\", \"41\\u00F76=\"],\n  [\"57\\u00F79=\", \"87\\u00F77=\"],\n  [\"82\\u00F78=\", \"67\\u00F73=\"],\n  [\"19\\u00F78=\", \"78\\u00F79=\"],\n  [\"94\\u00F74=\", \"23\\u00F73=\"],\n  [\"65\\u00F79=\", \"68\\u00F72=\"],\n  [\"74\\u00F72=\", \"48\\u00F79=\"],\n  [\"54\\u00F79=\", \"38\\u00F76=\"],\n];\n\n// 1) Kick off every search against the (still unmodified) document.\nconst searchResults = replacements.map(([oldText]) =>\n  context.document.body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// 2) Now that every match has been located, apply the replacements.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  if (items.length === 0) {\n    continue;\n  }\n  items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem cells in the table.\n# Cells are targeted by their (row, column) position in the single table\n# so the edit is unambiguous, even though some new values collide with\n# other old values used elsewhere in the table (e.g. 28\u00f75= is both an\n# old value at one cell and the new value written into another cell).\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph (first paragraph in the document).\n$d.Paragraphs(1).Range.Text = \"2025-09-19 Friday\"\n\n# 2) Table of division problems: 5 content rows (table rows 1, 5, 9, 13,\n# 17 \u2014 interleaved with blank spacer rows), 5 columns each.\n$tbl = $d.Tables(1)\n\n$rowValues = @{\n    1  = @(\"32\u00f74=\", \"61\u00f75=\", \"23\u00f77=\", \"80\u00f75=\", \"84\u00f74=\")\n    5  = @(\"74\u00f73=\", \"38\u00f79=\", \"87\u00f78=\", \"46\u00f76=\", \"20\u00f75=\")\n    9  = @(\"17\u00f79=\", \"16\u00f74=\", \"69\u00f74=\", \"37\u00f72=\", \"21\u00f74=\")\n    13 = @(\"28\u00f75=\", \"64\u00f75=\", \"41\u00f76=\", \"87\u00f77=\", \"67\u00f73=\")\n    17 = @(\"78\u00f79=\", \"23\u00f73=\", \"68\u00f72=\", \"48\u00f79=\", \"38\u00f76=\")\n}\n\nforeach ($rowIndex in $rowValues.Keys) {\n    $values = $rowValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $tbl.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
